# Produce an Excel sheet of (fake) "Wills" contact-data headers.
#
# The sheet currently has no data (dimension A1:A1, empty sheetData).
# We write a single header row spanning B1:J1, then apply a bold /
# bordered / center-top-aligned format to that header row so it reads
# like a table header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels, column B through J.
$headers = @(
    "first_name",
    "last_name",
    "email",
    "gender",
    "ip_address",
    "car_VIN",
    "city",
    "credit_card",
    "street_address"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Style the header row: bold text, thin box border around every cell,
# centered horizontally and top-aligned vertically.
$headerRange = $ws.Range("B1:J1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
